$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("CM1").Value = "1117_1400_new"

$values = @(
0.997335668580562,
0.99733510779404499,
0.99733455889208,
0.99732664382428804,
0.99732606216827402,
0.99732542237979405,
0.997323225393007,
0.99732265341488802,
0.99732217775908805,
0.99731994014295,
0.99731947717325298,
0.99731895075822896,
0.99731716311690399,
0.99731668146973895,
0.99731618642019204,
0.99731434229669702,
0.99731358138111004,
0.99731308809246699,
0.99731073938312897,
0.99731027772928305,
0.99730963447052201,
0.99730734269445498,
0.99730694386223295,
0.99730649074718203,
0.99730538441975602,
0.99730494817168003,
0.99730454214607001,
0.99729878843175601,
0.997296880057988,
0.99729490030429402,
0.99729307713079796,
0.99729106000397705,
0.99728976546127002,
0.99728784124992798,
0.99728608996956403,
0.99728341653441299,
0.99728113109618999,
0.99728035059815601,
0.997279486963904,
0.99727923669654905,
0.99727900693308802,
0.99727814247688096,
0.99727737509604097,
0.99727624425313699,
0.997275704620602,
0.99727551356421096,
0.99727535860291405,
0.99727483514434001,
0.99727429740612095,
0.99727379211918399,
0.997272599758109,
0.99727250625571595,
0.9972724226352,
0.99727166536177703,
0.99727128798813303,
0.99727098895161004,
0.99727098364046995,
0.99727098249939805
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 91).Value = $values[$i]
}

$ws.Range("CM2").Select()
